$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Neg_Change")
$ws2 = $wb.Worksheets.Item("Pos_Change")

# --- Neg_Change sheet: replace rows 2-14 ---
$ws1.Cells.Item(2, 1).Value = "JSWENERGY"
$ws1.Cells.Item(2, 2).Value = 534
$ws1.Cells.Item(2, 3).Value = 544
$ws1.Cells.Item(2, 4).Value = 532.55
$ws1.Cells.Item(2, 5).Value = 541.7
$ws1.Cells.Item(2, 6).Value = 1425566
$ws1.Cells.Item(2, 7).Value = 3256894
$ws1.Cells.Item(2, 8).Value = -0.5622927857031884
$ws1.Cells.Item(2, 9).Value = "JSWENERGY"
$ws1.Cells.Item(3, 1).Value = "INDHOTEL"
$ws1.Cells.Item(3, 2).Value = 733
$ws1.Cells.Item(3, 3).Value = 735
$ws1.Cells.Item(3, 4).Value = 728
$ws1.Cells.Item(3, 5).Value = 734.9
$ws1.Cells.Item(3, 6).Value = 1123403
$ws1.Cells.Item(3, 7).Value = 2734309
$ws1.Cells.Item(3, 8).Value = -0.5891455574333405
$ws1.Cells.Item(3, 9).Value = "INDHOTEL"
$ws1.Cells.Item(4, 1).Value = "IRFC"
$ws1.Cells.Item(4, 2).Value = 124.91
$ws1.Cells.Item(4, 3).Value = 125.35
$ws1.Cells.Item(4, 4).Value = 123.81
$ws1.Cells.Item(4, 5).Value = 125.1
$ws1.Cells.Item(4, 6).Value = 5774638
$ws1.Cells.Item(4, 7).Value = 12875915
$ws1.Cells.Item(4, 8).Value = -0.5515163000066403
$ws1.Cells.Item(4, 9).Value = "IRFC"
$ws1.Cells.Item(5, 1).Value = "CANBK"
$ws1.Cells.Item(5, 2).Value = 126.49
$ws1.Cells.Item(5, 3).Value = 126.99
$ws1.Cells.Item(5, 4).Value = 125.51
$ws1.Cells.Item(5, 5).Value = 126.24
$ws1.Cells.Item(5, 6).Value = 16981691
$ws1.Cells.Item(5, 7).Value = 38012061
$ws1.Cells.Item(5, 8).Value = -0.5532551891885052
$ws1.Cells.Item(5, 9).Value = "CANBK"
$ws1.Cells.Item(6, 1).Value = "AMBUJACEM"
$ws1.Cells.Item(6, 2).Value = 566.4
$ws1.Cells.Item(6, 3).Value = 569.15
$ws1.Cells.Item(6, 4).Value = 562.05
$ws1.Cells.Item(6, 5).Value = 565.8
$ws1.Cells.Item(6, 6).Value = 816362
$ws1.Cells.Item(6, 7).Value = 1866324
$ws1.Cells.Item(6, 8).Value = -0.5625829170069077
$ws1.Cells.Item(6, 9).Value = "AMBUJACEM"
$ws1.Cells.Item(7, 1).Value = "ABCAPITAL"
$ws1.Cells.Item(7, 2).Value = 296.9
$ws1.Cells.Item(7, 3).Value = 299
$ws1.Cells.Item(7, 4).Value = 294.95
$ws1.Cells.Item(7, 5).Value = 298.1
$ws1.Cells.Item(7, 6).Value = 2420997
$ws1.Cells.Item(7, 7).Value = 5342127
$ws1.Cells.Item(7, 8).Value = -0.5468102873630671
$ws1.Cells.Item(7, 9).Value = "ABCAPITAL"
$ws1.Cells.Item(8, 1).Value = "MANKIND"
$ws1.Cells.Item(8, 2).Value = 2450
$ws1.Cells.Item(8, 3).Value = 2484.2
$ws1.Cells.Item(8, 4).Value = 2448.1
$ws1.Cells.Item(8, 5).Value = 2476.8
$ws1.Cells.Item(8, 6).Value = 135618
$ws1.Cells.Item(8, 7).Value = 285681
$ws1.Cells.Item(8, 8).Value = -0.5252816953175045
$ws1.Cells.Item(8, 9).Value = "MANKIND"
$ws1.Cells.Item(9, 1).Value = "COFORGE"
$ws1.Cells.Item(9, 2).Value = 1720
$ws1.Cells.Item(9, 3).Value = 1728.5
$ws1.Cells.Item(9, 4).Value = 1699
$ws1.Cells.Item(9, 5).Value = 1723.1
$ws1.Cells.Item(9, 6).Value = 1167579
$ws1.Cells.Item(9, 7).Value = 2733508
$ws1.Cells.Item(9, 8).Value = -0.5728642462359722
$ws1.Cells.Item(9, 9).Value = "COFORGE"
$ws1.Cells.Item(10, 1).Value = "SUPREMEIND"
$ws1.Cells.Item(10, 2).Value = 4186.7
$ws1.Cells.Item(10, 3).Value = 4198.2
$ws1.Cells.Item(10, 4).Value = 4112.6
$ws1.Cells.Item(10, 5).Value = 4173
$ws1.Cells.Item(10, 6).Value = 52464
$ws1.Cells.Item(10, 7).Value = 116555
$ws1.Cells.Item(10, 8).Value = -0.5498777401226889
$ws1.Cells.Item(10, 9).Value = "SUPREMEIND"
$ws1.Cells.Item(11, 1).Value = "ICICIPRULI"
$ws1.Cells.Item(11, 2).Value = 592.1
$ws1.Cells.Item(11, 3).Value = 596.55
$ws1.Cells.Item(11, 4).Value = 589.5
$ws1.Cells.Item(11, 5).Value = 593.95
$ws1.Cells.Item(11, 6).Value = 461949
$ws1.Cells.Item(11, 7).Value = 972845
$ws1.Cells.Item(11, 8).Value = -0.5251566282398532
$ws1.Cells.Item(11, 9).Value = "ICICIPRULI"
$ws1.Cells.Item(12, 1).Value = "RVNL"
$ws1.Cells.Item(12, 2).Value = 344.5
$ws1.Cells.Item(12, 3).Value = 347.3
$ws1.Cells.Item(12, 4).Value = 341.2
$ws1.Cells.Item(12, 5).Value = 344.5
$ws1.Cells.Item(12, 6).Value = 4026015
$ws1.Cells.Item(12, 7).Value = 8765182
$ws1.Cells.Item(12, 8).Value = -0.5406809579082328
$ws1.Cells.Item(12, 9).Value = "RVNL"
$ws1.Cells.Item(13, 1).Value = "KFINTECH"
$ws1.Cells.Item(13, 2).Value = 1040
$ws1.Cells.Item(13, 3).Value = 1055.7
$ws1.Cells.Item(13, 4).Value = 1030
$ws1.Cells.Item(13, 5).Value = 1052.8
$ws1.Cells.Item(13, 6).Value = 565327
$ws1.Cells.Item(13, 7).Value = 1218429
$ws1.Cells.Item(13, 8).Value = -0.5360197434565329
$ws1.Cells.Item(13, 9).Value = "KFINTECH"
$ws1.Cells.Item(14, 1).Value = "POONAWALLA"
$ws1.Cells.Item(14, 2).Value = 526.55
$ws1.Cells.Item(14, 3).Value = 533.9
$ws1.Cells.Item(14, 4).Value = 524.2
$ws1.Cells.Item(14, 5).Value = 531
$ws1.Cells.Item(14, 6).Value = 1023794
$ws1.Cells.Item(14, 7).Value = 2234852
$ws1.Cells.Item(14, 8).Value = -0.5418962866444847
$ws1.Cells.Item(14, 9).Value = "POONAWALLA"

# --- Pos_Change sheet: replace rows 2-6 and add rows 7-15 ---
$ws2.Cells.Item(2, 1).Value = "HCLTECH"
$ws2.Cells.Item(2, 2).Value = 1453.2
$ws2.Cells.Item(2, 3).Value = 1491.1
$ws2.Cells.Item(2, 4).Value = 1452.5
$ws2.Cells.Item(2, 5).Value = 1484.5
$ws2.Cells.Item(2, 6).Value = 2681316
$ws2.Cells.Item(2, 7).Value = 1735807
$ws2.Cells.Item(2, 8).Value = 0.54470859951596
$ws2.Cells.Item(2, 9).Value = "HCLTECH"
$ws2.Cells.Item(3, 1).Value = "INFY"
$ws2.Cells.Item(3, 2).Value = 1502
$ws2.Cells.Item(3, 3).Value = 1514
$ws2.Cells.Item(3, 4).Value = 1487
$ws2.Cells.Item(3, 5).Value = 1509.7
$ws2.Cells.Item(3, 6).Value = 7028753
$ws2.Cells.Item(3, 7).Value = 4657192
$ws2.Cells.Item(3, 8).Value = 0.5092255161479278
$ws2.Cells.Item(3, 9).Value = "INFY"
$ws2.Cells.Item(4, 1).Value = "APOLLOHOSP"
$ws2.Cells.Item(4, 2).Value = 7690
$ws2.Cells.Item(4, 3).Value = 7750
$ws2.Cells.Item(4, 4).Value = 7620
$ws2.Cells.Item(4, 5).Value = 7700
$ws2.Cells.Item(4, 6).Value = 295802
$ws2.Cells.Item(4, 7).Value = 188810
$ws2.Cells.Item(4, 8).Value = 0.5666649012234521
$ws2.Cells.Item(4, 9).Value = "APOLLOHOSP"
$ws2.Cells.Item(5, 1).Value = "VEDL"
$ws2.Cells.Item(5, 2).Value = 473
$ws2.Cells.Item(5, 3).Value = 485
$ws2.Cells.Item(5, 4).Value = 467.35
$ws2.Cells.Item(5, 5).Value = 484.8
$ws2.Cells.Item(5, 6).Value = 16201074
$ws2.Cells.Item(5, 7).Value = 10158945
$ws2.Cells.Item(5, 8).Value = 0.5947594952034881
$ws2.Cells.Item(5, 9).Value = "VEDL"
$ws2.Cells.Item(6, 1).Value = "CGPOWER"
$ws2.Cells.Item(6, 2).Value = 743.9
$ws2.Cells.Item(6, 3).Value = 763.45
$ws2.Cells.Item(6, 4).Value = 741.8
$ws2.Cells.Item(6, 5).Value = 761.8
$ws2.Cells.Item(6, 6).Value = 2262502
$ws2.Cells.Item(6, 7).Value = 1563241
$ws2.Cells.Item(6, 8).Value = 0.4473149053792729
$ws2.Cells.Item(6, 9).Value = "CGPOWER"
$ws2.Cells.Item(7, 1).Value = "HAL"
$ws2.Cells.Item(7, 2).Value = 4792
$ws2.Cells.Item(7, 3).Value = 4863.2
$ws2.Cells.Item(7, 4).Value = 4744
$ws2.Cells.Item(7, 5).Value = 4857
$ws2.Cells.Item(7, 6).Value = 1003478
$ws2.Cells.Item(7, 7).Value = 704075
$ws2.Cells.Item(7, 8).Value = 0.425243049391045
$ws2.Cells.Item(7, 9).Value = "HAL"
$ws2.Cells.Item(8, 1).Value = "IOC"
$ws2.Cells.Item(8, 2).Value = 154.02
$ws2.Cells.Item(8, 3).Value = 157.2
$ws2.Cells.Item(8, 4).Value = 153.33
$ws2.Cells.Item(8, 5).Value = 155.1
$ws2.Cells.Item(8, 6).Value = 13999206
$ws2.Cells.Item(8, 7).Value = 9953844
$ws2.Cells.Item(8, 8).Value = 0.4064120353905486
$ws2.Cells.Item(8, 9).Value = "IOC"
$ws2.Cells.Item(9, 1).Value = "CHOLAFIN"
$ws2.Cells.Item(9, 2).Value = 1607.5
$ws2.Cells.Item(9, 3).Value = 1625
$ws2.Cells.Item(9, 4).Value = 1596.4
$ws2.Cells.Item(9, 5).Value = 1620.2
$ws2.Cells.Item(9, 6).Value = 1404177
$ws2.Cells.Item(9, 7).Value = 885008
$ws2.Cells.Item(9, 8).Value = 0.5866263355811473
$ws2.Cells.Item(9, 9).Value = "CHOLAFIN"
$ws2.Cells.Item(10, 1).Value = "ICICIGI"
$ws2.Cells.Item(10, 2).Value = 1880.2
$ws2.Cells.Item(10, 3).Value = 1894.8
$ws2.Cells.Item(10, 4).Value = 1854
$ws2.Cells.Item(10, 5).Value = 1890
$ws2.Cells.Item(10, 6).Value = 636816
$ws2.Cells.Item(10, 7).Value = 438884
$ws2.Cells.Item(10, 8).Value = 0.4509893274760529
$ws2.Cells.Item(10, 9).Value = "ICICIGI"
$ws2.Cells.Item(11, 1).Value = "NMDC"
$ws2.Cells.Item(11, 2).Value = 76.49
$ws2.Cells.Item(11, 3).Value = 79.24
$ws2.Cells.Item(11, 4).Value = 76.49
$ws2.Cells.Item(11, 5).Value = 78.82
$ws2.Cells.Item(11, 6).Value = 33869311
$ws2.Cells.Item(11, 7).Value = 21657726
$ws2.Cells.Item(11, 8).Value = 0.5638442835595944
$ws2.Cells.Item(11, 9).Value = "NMDC"
$ws2.Cells.Item(12, 1).Value = "LTF"
$ws2.Cells.Item(12, 2).Value = 255.95
$ws2.Cells.Item(12, 3).Value = 261
$ws2.Cells.Item(12, 4).Value = 255.55
$ws2.Cells.Item(12, 5).Value = 260.19
$ws2.Cells.Item(12, 6).Value = 3597676
$ws2.Cells.Item(12, 7).Value = 2485985
$ws2.Cells.Item(12, 8).Value = 0.4471833096338071
$ws2.Cells.Item(12, 9).Value = "LTF"
$ws2.Cells.Item(13, 1).Value = "PAGEIND"
$ws2.Cells.Item(13, 2).Value = 42000
$ws2.Cells.Item(13, 3).Value = 42085
$ws2.Cells.Item(13, 4).Value = 41365
$ws2.Cells.Item(13, 5).Value = 41700
$ws2.Cells.Item(13, 6).Value = 20885
$ws2.Cells.Item(13, 7).Value = 13956
$ws2.Cells.Item(13, 8).Value = 0.4964889653195758
$ws2.Cells.Item(13, 9).Value = "PAGEIND"
$ws2.Cells.Item(14, 1).Value = "LICHSGFIN"
$ws2.Cells.Item(14, 2).Value = 563.5
$ws2.Cells.Item(14, 3).Value = 565.55
$ws2.Cells.Item(14, 4).Value = 559
$ws2.Cells.Item(14, 5).Value = 560
$ws2.Cells.Item(14, 6).Value = 2247538
$ws2.Cells.Item(14, 7).Value = 1591967
$ws2.Cells.Item(14, 8).Value = 0.4117993651878462
$ws2.Cells.Item(14, 9).Value = "LICHSGFIN"
$ws2.Cells.Item(15, 1).Value = "DELHIVERY"
$ws2.Cells.Item(15, 2).Value = 467.3
$ws2.Cells.Item(15, 3).Value = 475.5
$ws2.Cells.Item(15, 4).Value = 466
$ws2.Cells.Item(15, 5).Value = 468.9
$ws2.Cells.Item(15, 6).Value = 3354443
$ws2.Cells.Item(15, 7).Value = 2249868
$ws2.Cells.Item(15, 8).Value = 0.4909510246823369
$ws2.Cells.Item(15, 9).Value = "DELHIVERY"

Write-Output "Updated Neg_Change rows 2-14 and Pos_Change rows 2-15"
